# "modified order in big_five data"
#
# The four colored row-blocks inside the data region (rows 38-72, columns
# A:E) of the big_five worksheet get reordered. Each block keeps its own
# fill/style as it moves:
#
#   block A = rows 38-53 (16 rows, "theme7" fill)
#   block B = rows 54-63 (10 rows, "green"  fill)
#   block C = rows 64-65 ( 2 rows, "theme9" fill)
#   block D = rows 66-72 ( 7 rows, "theme5" fill)
#
# Original order:  A, B, C, D
# New order:        C, A, D, B
#
# i.e. new layout starting at row 38:
#   C -> rows 38-39
#   A -> rows 40-55
#   D -> rows 56-62
#   B -> rows 63-72
#
# Range.Copy(destination) is used because it is the one operation in this
# COM surface that reliably carries both the values AND the existing cell
# style (fill) index across without duplicating style entries. Because the
# four blocks' source and destination ranges overlap each other once
# rearranged, every block is first staged to a scratch area well below the
# sheet's used range (rows 200+) and then copied back into its final spot,
# which avoids a source block being clobbered before it has been consumed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stage each block (values + style) out of the way ---------------------
$ws.Range("A38:E53").Copy($ws.Range("A201"))   # block A (16 rows) -> staged at 201-216
$ws.Range("A54:E63").Copy($ws.Range("A221"))   # block B (10 rows) -> staged at 221-230
$ws.Range("A64:E65").Copy($ws.Range("A241"))   # block C ( 2 rows) -> staged at 241-242
$ws.Range("A66:E72").Copy($ws.Range("A251"))   # block D ( 7 rows) -> staged at 251-257

# --- Copy staged blocks back into their new positions ----------------------
$ws.Range("A241:E242").Copy($ws.Range("A38"))  # block C -> rows 38-39
$ws.Range("A201:E216").Copy($ws.Range("A40"))  # block A -> rows 40-55
$ws.Range("A251:E257").Copy($ws.Range("A56"))  # block D -> rows 56-62
$ws.Range("A221:E230").Copy($ws.Range("A63"))  # block B -> rows 63-72

# --- Clean up the scratch/staging area --------------------------------------
$ws.Range("A201:E260").Clear()

# --- Update the view's selection to match the new scroll/selection state ---
$ws.Activate()
$ws.Range("A63:E72").Select()
